$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 874.7794
$ws.Range("J17").Value = 874.7794
$ws.Range("L17").Value = 2624.3382
$ws.Range("N17").Value = -2960.3382
$ws.Range("H86").Value = 1610.2222
$ws.Range("I86").Value = 1498.8572
$ws.Range("K86").Value = 1498.8572
$ws.Range("M86").Value = -375.8571999999999
$ws.Range("H89").Value = 1610.2222
$ws.Range("I89").Value = 1498.8572
$ws.Range("K89").Value = 7494.286
$ws.Range("M89").Value = -1878.286
$ws.Range("H137").Value = 37254.07
$ws.Range("I137").Value = 1448.1428
$ws.Range("J137").Value = 144671.86
$ws.Range("K137").Value = 4344.428400000001
$ws.Range("L137").Value = 434015.58
$ws.Range("M137").Value = -1794.428400000001
$ws.Range("N137").Value = -439115.58
$ws.Range("H138").Value = 3376.3872
$ws.Range("I138").Value = 3103.8696
$ws.Range("J138").Value = 4159.875
$ws.Range("K138").Value = 9311.6088
$ws.Range("L138").Value = 12479.625
$ws.Range("M138").Value = -4171.6088
$ws.Range("N138").Value = -22759.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2908051
$ws.Range("J2").Value = 1333
$ws.Range("L2").Value = 1333
$ws.Range("N2").Value = -1559
$ws.Range("H32").Value = 2688.2273
$ws.Range("I32").Value = 2046.9865
$ws.Range("K32").Value = 2046.9865
$ws.Range("M32").Value = -1759.9865
$ws.Range("H45").Value = 1653.7333
$ws.Range("I45").Value = 1112.25
$ws.Range("K45").Value = 1112.25
$ws.Range("M45").Value = -735.25
$ws.Range("H74").Value = 1394.742
$ws.Range("I74").Value = 1151.2273
$ws.Range("J74").Value = 1990
$ws.Range("K74").Value = 1151.2273
$ws.Range("L74").Value = 1990
$ws.Range("M74").Value = -277.2273
$ws.Range("N74").Value = -3738
$ws.Range("H77").Value = 1394.742
$ws.Range("I77").Value = 1151.2273
$ws.Range("J77").Value = 1990
$ws.Range("K77").Value = 5756.136500000001
$ws.Range("L77").Value = 9950
$ws.Range("M77").Value = -1388.136500000001
$ws.Range("N77").Value = -18686
$ws.Range("H110").Value = 2953.6667
$ws.Range("J110").Value = 5505.5
$ws.Range("L110").Value = 5505.5
$ws.Range("N110").Value = -9595.5
$ws.Range("H116").Value = 2908051
$ws.Range("J116").Value = 1333
$ws.Range("L116").Value = 1333
$ws.Range("N116").Value = -5921
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2908051
$ws.Range("J3").Value = 1333
$ws.Range("L3").Value = 1333
$ws.Range("N3").Value = -1561
$ws.Range("H94").Value = 582.6429000000001
$ws.Range("I94").Value = 419.76
$ws.Range("K94").Value = 419.76
$ws.Range("M94").Value = 31.24000000000001
$ws.Range("H105").Value = 2187.423
$ws.Range("I105").Value = 2140.5833
$ws.Range("K105").Value = 2140.5833
$ws.Range("M105").Value = -393.5832999999998
$ws.Range("H107").Value = 1746.95
$ws.Range("I107").Value = 1676
$ws.Range("J107").Value = 1959.8
$ws.Range("K107").Value = 1676
$ws.Range("L107").Value = 1959.8
$ws.Range("M107").Value = 244
$ws.Range("N107").Value = -5799.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 7500
$ws.Range("J29").Value = 7500
$ws.Range("L29").Value = 7500
$ws.Range("N29").Value = -8086
$ws.Range("H31").Value = 1902.6471
$ws.Range("I31").Value = 1477.1818
$ws.Range("K31").Value = 1477.1818
$ws.Range("M31").Value = -1182.1818
$ws.Range("H34").Value = 1902.6471
$ws.Range("I34").Value = 1477.1818
$ws.Range("K34").Value = 1477.1818
$ws.Range("M34").Value = -1275.1818
$ws.Range("H58").Value = 1612279.6
$ws.Range("I58").Value = 2289759.5
$ws.Range("K58").Value = 2289759.5
$ws.Range("M58").Value = -2289556.5
$ws.Range("H107").Value = 419.0357
$ws.Range("I107").Value = 352.2857
$ws.Range("J107").Value = 619.2857
$ws.Range("K107").Value = 352.2857
$ws.Range("L107").Value = 619.2857
$ws.Range("M107").Value = 1567.7143
$ws.Range("N107").Value = -4459.2857
$ws.Range("H132").Value = 2463.8147
$ws.Range("I132").Value = 1462.3846
$ws.Range("K132").Value = 4387.1538
$ws.Range("M132").Value = -1857.1538
$ws.Range("H136").Value = 1612279.6
$ws.Range("I136").Value = 2289759.5
$ws.Range("K136").Value = 6869278.5
$ws.Range("M136").Value = -6866728.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 574.7778
$ws.Range("I26").Value = 623.8
$ws.Range("J26").Value = 513.5
$ws.Range("K26").Value = 1871.4
$ws.Range("L26").Value = 1540.5
$ws.Range("M26").Value = -1583.4
$ws.Range("N26").Value = -2116.5
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H107").Value = 658.3333
$ws.Range("I107").Value = 265.57144
$ws.Range("J107").Value = 820.05884
$ws.Range("K107").Value = 796.71432
$ws.Range("L107").Value = 2460.17652
$ws.Range("M107").Value = 1123.28568
$ws.Range("N107").Value = -6300.17652
$ws.Range("H131").Value = 807.09186
$ws.Range("J131").Value = 810.36456
$ws.Range("L131").Value = 2431.09368
$ws.Range("N131").Value = -12511.09368
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13069.111
$ws.Range("I70").Value = 19824.4
$ws.Range("J70").Value = 4625
$ws.Range("K70").Value = 19824.4
$ws.Range("L70").Value = 4625
$ws.Range("M70").Value = -19554.4
$ws.Range("N70").Value = -5165
$ws.Range("H73").Value = 13069.111
$ws.Range("I73").Value = 19824.4
$ws.Range("J73").Value = 4625
$ws.Range("K73").Value = 19824.4
$ws.Range("L73").Value = 4625
$ws.Range("M73").Value = -18888.4
$ws.Range("N73").Value = -6497
$ws.Range("H80").Value = 3932.6
$ws.Range("I80").Value = 3916
$ws.Range("J80").Value = 3999
$ws.Range("K80").Value = 3916
$ws.Range("L80").Value = 3999
$ws.Range("M80").Value = -2918
$ws.Range("N80").Value = -5995
$ws.Range("H83").Value = 3932.6
$ws.Range("I83").Value = 3916
$ws.Range("J83").Value = 3999
$ws.Range("K83").Value = 19580
$ws.Range("L83").Value = 19995
$ws.Range("M83").Value = -14588
$ws.Range("N83").Value = -29979
$ws.Range("H93").Value = 29624.5
$ws.Range("J93").Value = 29624.5
$ws.Range("L93").Value = 29624.5
$ws.Range("N93").Value = -33368.5
$ws.Range("H132").Value = 1167940.8
$ws.Range("I132").Value = 1750150
$ws.Range("J132").Value = 3522.4546
$ws.Range("K132").Value = 5250450
$ws.Range("L132").Value = 10567.3638
$ws.Range("M132").Value = -5247920
$ws.Range("N132").Value = -15627.3638
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 6195.091
$ws.Range("I4").Value = 5314.6
$ws.Range("J4").Value = 15000
$ws.Range("K4").Value = 5314.6
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = -5201.6
$ws.Range("N4").Value = -15226
$ws.Range("H28").Value = 6195.091
$ws.Range("I28").Value = 5314.6
$ws.Range("J28").Value = 15000
$ws.Range("K28").Value = 5314.6
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = -5082.6
$ws.Range("N28").Value = -15464
$ws.Range("H37").Value = 6195.091
$ws.Range("I37").Value = 5314.6
$ws.Range("J37").Value = 15000
$ws.Range("K37").Value = 5314.6
$ws.Range("L37").Value = 15000
$ws.Range("M37").Value = -5207.6
$ws.Range("N37").Value = -15214
$ws.Range("H40").Value = 4464.7334
$ws.Range("I40").Value = 1864
$ws.Range("J40").Value = 8365.833000000001
$ws.Range("K40").Value = 1864
$ws.Range("L40").Value = 8365.833000000001
$ws.Range("M40").Value = -1728
$ws.Range("N40").Value = -8637.833000000001
$ws.Range("H68").Value = 3082.5833
$ws.Range("I68").Value = 2621.2222
$ws.Range("J68").Value = 4466.6665
$ws.Range("K68").Value = 2621.2222
$ws.Range("L68").Value = 4466.6665
$ws.Range("M68").Value = -1872.2222
$ws.Range("N68").Value = -5964.6665
$ws.Range("H71").Value = 3082.5833
$ws.Range("I71").Value = 2621.2222
$ws.Range("J71").Value = 4466.6665
$ws.Range("K71").Value = 13106.111
$ws.Range("L71").Value = 22333.3325
$ws.Range("M71").Value = -9362.111000000001
$ws.Range("N71").Value = -29821.3325
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 20964.928
$ws.Range("J2").Value = 20964.928
$ws.Range("L2").Value = 20964.928
$ws.Range("N2").Value = -21188.928
$ws.Range("H30").Value = 11350
$ws.Range("I30").Value = 500
$ws.Range("J30").Value = 14966.667
$ws.Range("K30").Value = 500
$ws.Range("L30").Value = 14966.667
$ws.Range("M30").Value = -393
$ws.Range("N30").Value = -15180.667
$ws.Range("H126").Value = 14700.4
$ws.Range("I126").Value = 22800.8
$ws.Range("J126").Value = 6600
$ws.Range("K126").Value = 68402.39999999999
$ws.Range("L126").Value = 19800
$ws.Range("M126").Value = -65932.39999999999
$ws.Range("N126").Value = -24740
$ws.Range("H132").Value = 1115.1316
$ws.Range("J132").Value = 2114.8572
$ws.Range("L132").Value = 6344.571599999999
$ws.Range("N132").Value = -11404.5716

Write-Output "Applied 235 cell updates across 8 sheets"